$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 991.0769
$ws.Range("I40").Value = 1183.7142
$ws.Range("J40").Value = 766.3333
$ws.Range("K40").Value = 1183.7142
$ws.Range("L40").Value = 766.3333
$ws.Range("M40").Value = -1008.7142
$ws.Range("N40").Value = -1116.3333

$ws.Range("H74").Value = 4225
$ws.Range("I74").Value = 3633.3333
$ws.Range("J74").Value = 6000
$ws.Range("K74").Value = 3633.3333
$ws.Range("L74").Value = 6000
$ws.Range("M74").Value = -2697.3333
$ws.Range("N74").Value = -7872

$ws.Range("H77").Value = 4225
$ws.Range("I77").Value = 3633.3333
$ws.Range("J77").Value = 6000
$ws.Range("K77").Value = 18166.6665
$ws.Range("L77").Value = 30000
$ws.Range("M77").Value = -13486.6665
$ws.Range("N77").Value = -39360

$ws.Range("H94").Value = 3330.1667
$ws.Range("I94").Value = 2996.2
$ws.Range("J94").Value = 5000
$ws.Range("K94").Value = 2996.2
$ws.Range("L94").Value = 5000
$ws.Range("M94").Value = -2545.2
$ws.Range("N94").Value = -5902

$ws.Range("H98").Value = 3126.2727
$ws.Range("I98").Value = 3126.2727
$ws.Range("K98").Value = 3126.2727
$ws.Range("M98").Value = -1628.2727

$ws.Range("H122").Value = 3126.2727
$ws.Range("I122").Value = 3126.2727
$ws.Range("K122").Value = 9378.8181
$ws.Range("M122").Value = -6928.8181

$ws.Range("H132").Value = 1038.7715
$ws.Range("I132").Value = 974.4545000000001
$ws.Range("J132").Value = 2100
$ws.Range("K132").Value = 2923.3635
$ws.Range("L132").Value = 6300
$ws.Range("M132").Value = -393.3635000000004
$ws.Range("N132").Value = -11360

$ws.Range("H134").Value = 48171.8
$ws.Range("J134").Value = 48171.8
$ws.Range("L134").Value = 48171.8
$ws.Range("N134").Value = -58311.8

$ws.Range("H137").Value = 2499.2354
$ws.Range("I137").Value = 1403.8334
$ws.Range("K137").Value = 4211.5002
$ws.Range("M137").Value = -1661.5002

$ws.Range("H138").Value = 3417.8333
$ws.Range("I138").Value = 3765.9412
$ws.Range("J138").Value = 3106.3684
$ws.Range("K138").Value = 11297.8236
$ws.Range("L138").Value = 9319.1052
$ws.Range("M138").Value = -6157.8236
$ws.Range("N138").Value = -19599.1052

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3094.4727
$ws.Range("I32").Value = 2848.3044
$ws.Range("K32").Value = 2848.3044
$ws.Range("M32").Value = -2561.3044

$ws.Range("H74").Value = 1719.9375
$ws.Range("I74").Value = 1553.2
$ws.Range("J74").Value = 1997.8334
$ws.Range("K74").Value = 1553.2
$ws.Range("L74").Value = 1997.8334
$ws.Range("M74").Value = -679.2
$ws.Range("N74").Value = -3745.8334

$ws.Range("H77").Value = 1719.9375
$ws.Range("I77").Value = 1553.2
$ws.Range("J77").Value = 1997.8334
$ws.Range("K77").Value = 7766
$ws.Range("L77").Value = 9989.166999999999
$ws.Range("M77").Value = -3398
$ws.Range("N77").Value = -18725.167

$ws.Range("H102").Value = 1557.5
$ws.Range("I102").Value = 1000
$ws.Range("K102").Value = 1000
$ws.Range("M102").Value = 622

$ws.Range("H123").Value = 63998.5
$ws.Range("J123").Value = 63998.5
$ws.Range("L123").Value = 63998.5
$ws.Range("N123").Value = -73798.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2320.261
$ws.Range("I105").Value = 2356.9412
$ws.Range("J105").Value = 2216.3333
$ws.Range("K105").Value = 2356.9412
$ws.Range("L105").Value = 2216.3333
$ws.Range("M105").Value = -609.9412000000002
$ws.Range("N105").Value = -5710.3333

$ws.Range("H107").Value = 1953.2
$ws.Range("I107").Value = 1679.9
$ws.Range("K107").Value = 1679.9
$ws.Range("M107").Value = 240.0999999999999

$ws.Range("H134").Value = 7305.926
$ws.Range("I134").Value = 8668.1
$ws.Range("J134").Value = 3414
$ws.Range("K134").Value = 26004.3
$ws.Range("L134").Value = 10242
$ws.Range("M134").Value = -23469.3
$ws.Range("N134").Value = -15312

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 4349590
$ws.Range("I58").Value = 4349590
$ws.Range("K58").Value = 4349590
$ws.Range("M58").Value = -4349387

$ws.Range("H99").Value = 2293.5833
$ws.Range("J99").Value = 2333.3333
$ws.Range("L99").Value = 2333.3333
$ws.Range("N99").Value = -5329.3333

$ws.Range("H105").Value = 1354.3334
$ws.Range("I105").Value = 1170.5714
$ws.Range("K105").Value = 1170.5714
$ws.Range("M105").Value = 576.4286

$ws.Range("H126").Value = 2293.5833
$ws.Range("J126").Value = 2333.3333
$ws.Range("L126").Value = 6999.999899999999
$ws.Range("N126").Value = -11939.9999

$ws.Range("H132").Value = 1892.4783
$ws.Range("I132").Value = 1172.5625
$ws.Range("K132").Value = 3517.6875
$ws.Range("M132").Value = -987.6875

$ws.Range("H136").Value = 4349590
$ws.Range("I136").Value = 4349590
$ws.Range("K136").Value = 13048770
$ws.Range("M136").Value = -13046220

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 6611.9414
$ws.Range("I34").Value = 16873.166
$ws.Range("J34").Value = 1014.9091
$ws.Range("K34").Value = 50619.49800000001
$ws.Range("L34").Value = 3044.7273
$ws.Range("M34").Value = -50535.49800000001
$ws.Range("N34").Value = -3212.7273

$ws.Range("H48").Value = 2449.25
$ws.Range("J48").Value = 2449.25
$ws.Range("L48").Value = 7347.75
$ws.Range("N48").Value = -7847.75

$ws.Range("H107").Value = 864.9375
$ws.Range("J107").Value = 1067.4445
$ws.Range("L107").Value = 3202.3335
$ws.Range("N107").Value = -7042.333500000001

$ws.Range("H127").Value = 1741.5
$ws.Range("J127").Value = 1741.5
$ws.Range("L127").Value = 5224.5
$ws.Range("N127").Value = -15144.5

$ws.Range("H132").Value = 1251.8334
$ws.Range("J132").Value = 1251.8334
$ws.Range("L132").Value = 11266.5006
$ws.Range("N132").Value = -16326.5006

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1056.625
$ws.Range("J113").Value = 1203.25
$ws.Range("L113").Value = 1203.25
$ws.Range("N113").Value = -5543.25

$ws.Range("H122").Value = 1492.8182
$ws.Range("I122").Value = 1062.1428
$ws.Range("J122").Value = 2246.5
$ws.Range("K122").Value = 3186.4284
$ws.Range("L122").Value = 6739.5
$ws.Range("M122").Value = -736.4284000000002
$ws.Range("N122").Value = -11639.5

$ws.Range("H126").Value = 3537316.8
$ws.Range("I126").Value = 13892392
$ws.Range("K126").Value = 41677176
$ws.Range("M126").Value = -41674706

$ws.Range("H127").Value = 33867.555
$ws.Range("J127").Value = 33867.555
$ws.Range("L127").Value = 33867.555
$ws.Range("N127").Value = -43787.555

$ws.Range("H132").Value = 1605286.4
$ws.Range("I132").Value = 2264383.8
$ws.Range("K132").Value = 6793151.399999999
$ws.Range("M132").Value = -6790621.399999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1727.9565
$ws.Range("I61").Value = 1705.1333
$ws.Range("J61").Value = 1770.75
$ws.Range("K61").Value = 1705.1333
$ws.Range("L61").Value = 1770.75
$ws.Range("M61").Value = -1503.1333
$ws.Range("N61").Value = -2174.75

$ws.Range("H113").Value = 1727.9565
$ws.Range("I113").Value = 1705.1333
$ws.Range("J113").Value = 1770.75
$ws.Range("K113").Value = 1705.1333
$ws.Range("L113").Value = 1770.75
$ws.Range("M113").Value = 464.8667
$ws.Range("N113").Value = -6110.75

$ws.Range("H136").Value = 2662.182
$ws.Range("I136").Value = 1447.04
$ws.Range("K136").Value = 4341.12
$ws.Range("M136").Value = -1791.12

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H58").Value = 0
$ws.Range("I58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("M58").ClearContents()

$ws.Range("H81").Value = 3000
$ws.Range("I81").Value = 3000
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 6000
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -4939
$ws.Range("N81").ClearContents()

$ws.Range("H84").Value = 3000
$ws.Range("I84").Value = 3000
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 30000
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -24696
$ws.Range("N84").ClearContents()

$ws.Range("H122").Value = 56697.145
$ws.Range("I122").Value = 78936
$ws.Range("K122").Value = 236808
$ws.Range("M122").Value = -234358

$ws.Range("H126").Value = 4463.8
$ws.Range("I126").Value = 2781.6
$ws.Range("K126").Value = 8344.799999999999
$ws.Range("M126").Value = -5874.799999999999

$ws.Range("H132").Value = 1963.8125
$ws.Range("I132").Value = 1115.0952
$ws.Range("J132").Value = 3584.0908
$ws.Range("K132").Value = 3345.2856
$ws.Range("L132").Value = 10752.2724
$ws.Range("M132").Value = -815.2856000000002
$ws.Range("N132").Value = -15812.2724
